$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old contents (values only, keep formatting) for the full previous range
$ws.Range("A1:B34").ClearContents()

# Re-write header row
$ws.Range("A1").Value = "Cluster name"
$ws.Range("B1").Value = "Active cases"

# Re-write data rows in the new (alphabetical) cluster order
$ws.Range("A2").Value = '126 Racecourse Road Public Housing Tower Flemington'
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = '139 Highett St Apartment Complex Richmond'
$ws.Range("B3").Value = 10
$ws.Range("A4").Value = '3175 The Bays Aged Care Facility Hastings'
$ws.Range("B4").Value = 17
$ws.Range("A5").Value = '3612 BlueCross Glengowrie Outbreak'
$ws.Range("B5").Value = 10
$ws.Range("A6").Value = '3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak'
$ws.Range("B6").Value = 20
$ws.Range("A7").Value = '4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak'
$ws.Range("B7").Value = 14
$ws.Range("A8").Value = 'Australian Lamb Colac East'
$ws.Range("B8").Value = 14
$ws.Range("A9").Value = 'Bread Solutions Braeside Outbreak'
$ws.Range("B9").Value = 20
$ws.Range("A10").Value = 'CS Square Caroline Springs'
$ws.Range("B10").Value = 15
$ws.Range("A11").Value = 'Carton Finishing Pty. Ltd. Campbellfield'
$ws.Range("B11").Value = 10
$ws.Range("A12").Value = 'Cedar Meats Australia Brooklyn Outbreak'
$ws.Range("B12").Value = 11
$ws.Range("A13").Value = 'Community Kids Bayswater Early Education Centre Bayswater North'
$ws.Range("B13").Value = 10
$ws.Range("A14").Value = 'Embracia Aged Care Reservoir Outbreak'
$ws.Range("B14").Value = 23
$ws.Range("A15").Value = 'FedEx Station Melbourne Airport'
$ws.Range("B15").Value = 11
$ws.Range("A16").Value = 'Green Leaves Early Learning Centre Highlands Craigieburn'
$ws.Range("B16").Value = 10
$ws.Range("A17").Value = 'Guardian Childcare Caulfield Outbreak'
$ws.Range("B17").Value = 15
$ws.Range("A18").Value = 'Kool Kidz Childcare Narre Warren'
$ws.Range("B18").Value = 11
$ws.Range("A19").Value = 'Lantmannen Unibake Australia Mordialloc'
$ws.Range("B19").Value = 26
$ws.Range("A20").Value = 'Nido Early School Ascot Vale'
$ws.Range("B20").Value = 15
$ws.Range("A21").Value = 'Nido Early School Glenroy'
$ws.Range("B21").Value = 23
$ws.Range("A22").Value = 'Northern Health Northern Hospital Epping Emergency Department Tier 1B'
$ws.Range("B22").Value = 48
$ws.Range("A23").Value = 'Northern Health The Northern Hospital Epping'
$ws.Range("B23").Value = 16
$ws.Range("A24").Value = 'Oceania Meat Processors Laverton North Outbreak'
$ws.Range("B24").Value = 13
$ws.Range("A25").Value = 'Pick It Up Fitness Mulgrave Outbreak'
$ws.Range("B25").Value = 11
$ws.Range("A26").Value = 'Social Gathering Warrnambool 28 Sep Outbreak'
$ws.Range("B26").Value = 17
$ws.Range("A27").Value = 'St Vincents Hospital Emergency Department Melbourne'
$ws.Range("B27").Value = 37
$ws.Range("A28").Value = 'Target Distribution Centre Truganina Outbreak'
$ws.Range("B28").Value = 21
$ws.Range("A29").Value = 'The Robin Hood Inn Drouin West Outbreak'
$ws.Range("B29").Value = 16
$ws.Range("A30").Value = 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'
$ws.Range("B30").Value = 13
$ws.Range("A31").Value = 'Visy Recycling Springvale'
$ws.Range("B31").Value = 31
$ws.Range("A32").Value = 'Werribee Mercy Hospital Emergency Department'
$ws.Range("B32").Value = 22
$ws.Range("A33").Value = 'Western Health Sunshine Hospital Emergency Department'
$ws.Range("B33").Value = 18
